$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.820.97'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '1.892.71'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.12'
$ws.Range('E5').Value = '  +2.42%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4812'
$ws.Range('E7').Value = '  +2.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2966'
$ws.Range('E8').Value = '  +7.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06660'
$ws.Range('E9').Value = '  +4.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.65'
$ws.Range('E10').Value = '  +6.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '100.28'
$ws.Range('E11').Value = '  +17.75%  '
$ws.Range('D12').Value = '1.874.31'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('E13').Value = '  +1.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.171'
$ws.Range('E14').Value = '  +3.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6606'
$ws.Range('E15').Value = '  +4.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '305.87'
$ws.Range('E16').Value = '  +26.62%  '
$ws.Range('D17').Value = '30.796.17'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('E18').Value = '  +3.75%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007642'
$ws.Range('E19').Value = '  +3.55%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').Value = '2.117.67'
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.180'
$ws.Range('E23').Value = '  +3.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.241'
$ws.Range('E24').Value = '  +4.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.351'
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.83'
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.36'
$ws.Range('E27').Value = '  +12.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1151'
$ws.Range('E28').Value = '  +11.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.961'
$ws.Range('E29').Value = '  +3.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.354'
$ws.Range('E30').Value = '  -1.65%  '
$ws.Range('E31').Value = '  +2.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.019'
$ws.Range('E32').Value = '  +4.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05097'
$ws.Range('E33').Value = '  +3.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7458'
$ws.Range('E34').Value = '  +5.14%  '
$ws.Range('E35').Value = '  +1.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.709'
$ws.Range('E36').Value = '  +0.44%  '
$ws.Range('E37').Value = '  +3.21%  '
$ws.Range('E38').Value = '  +0.76%  '
$ws.Range('E39').Value = '  +3.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8971'
$ws.Range('E40').Value = '  +1.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '108.34'
$ws.Range('E41').Value = '  +2.19%  '
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4223'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.667'
$ws.Range('E44').Value = '  +1.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '68.09'
$ws.Range('E45').Value = '  +9.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.402'
$ws.Range('E46').Value = '  +2.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.141'
$ws.Range('E47').Value = '  +5.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1237'
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.13'
$ws.Range('E49').Value = '  +4.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05652'
$ws.Range('E50').Value = '  +1.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.404'
$ws.Range('E51').Value = '  +1.51%  '
